$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 111
$ws1.Range("H2").Value = 5.48
$ws1.Range("L2").Value = 0.83

# Row 3
$ws1.Range("D3").Value = 115
$ws1.Range("H3").Value = 4.33
$ws1.Range("L3").Value = 0.83

# Row 4
$ws1.Range("D4").Value = 118
$ws1.Range("H4").Value = 3.26
$ws1.Range("L4").Value = 1.1

# Row 5
$ws1.Range("D5").Value = 115
$ws1.Range("H5").Value = 2.31
$ws1.Range("L5").Value = 0.87

# Row 6
$ws1.Range("D6").Value = 117
$ws1.Range("H6").Value = 1.29
$ws1.Range("L6").Value = 0.8100000000000001

# Row 7
$ws1.Range("D7").Value = 113
$ws1.Range("H7").Value = 0.31
$ws1.Range("I7").Value = "High"
$ws1.Range("L7").Value = 1.02

# Row 8
$ws1.Range("D8").Value = 118
$ws1.Range("L8").Value = 0.95

# Row 9
$ws1.Range("L9").Value = 1.02

# Row 10
$ws1.Range("L10").Value = 0.97

# Row 11
$ws1.Range("L11").Value = 0.89

# Row 12
$ws1.Range("L12").Value = 0.88

# Row 13
$ws1.Range("L13").Value = 1.14

# Row 15
$ws1.Range("L15").Value = 0.82

# Row 16
$ws1.Range("L16").Value = 0.9399999999999999

# Row 17
$ws1.Range("L17").Value = 1.11

# --- Sheet: Summary ---
# These cells hold their numbers as text (inline strings) in the workbook,
# so a leading apostrophe is used to force Excel to keep them as text
# instead of auto-converting the numeric-looking string to a number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'1808"
$ws2.Range("B10").Value = "'928"
$ws2.Range("B11").Value = "'462"
$ws2.Range("B14").Value = "'104"
